# KHL stats refresh: 2025-11-27 matches added, SOG aggregates recomputed,
# Meta_ext as_of_utc bumped. Mirrors the daily "publish files + archive" job.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Matches_SOG: fix F324 (sog_away 24 -> 25) and append 4 new matches
#    played 2025-11-27 (rows 326-329).
# ---------------------------------------------------------------------------
$wsMatches = $wb.Worksheets.Item("Matches_SOG")

$wsMatches.Cells.Item(324, 6).Value = 25

$newMatches = @(
  @{row=326; uid="897819"; date="2025-11-27T17:00:00"; home="Салават Юлаев"; away="Барыс";         sogH=35; sogA=33},
  @{row=327; uid="897821"; date="2025-11-27T17:00:00"; home="Металлург Мг";  away="Авангард";       sogH=60; sogA=30},
  @{row=328; uid="897820"; date="2025-11-27T19:30:00"; home="Спартак";       away="Автомобилист";   sogH=25; sogA=35},
  @{row=329; uid="897822"; date="2025-11-27T19:30:00"; home="ЦСКА";         away="Лада";            sogH=39; sogA=19}
)

foreach ($m in $newMatches) {
  $r = $m.row
  # uid looks like a plain integer; force text formatting first so the
  # COM layer doesn't silently coerce it to a Double (matches the source
  # file, where uid is stored as a string column).
  $wsMatches.Cells.Item($r, 1).NumberFormat = "@"
  $wsMatches.Cells.Item($r, 1).Value = $m.uid
  $wsMatches.Cells.Item($r, 2).Value = $m.date
  $wsMatches.Cells.Item($r, 3).Value = $m.home
  $wsMatches.Cells.Item($r, 4).Value = $m.away
  $wsMatches.Cells.Item($r, 5).Value = $m.sogH
  $wsMatches.Cells.Item($r, 6).Value = $m.sogA
  $wsMatches.Cells.Item($r, 7).Value = "khl_text"
}

# ---------------------------------------------------------------------------
# 2) Shots_HA: as_of_utc bumped for every team (2025-11-26 -> 2025-11-27),
#    plus the home/away shot tallies for teams who played on the 27th.
#    Columns: A Team, C season_id, D as_of_utc,
#             E GP_home, F GP_away, G HOGF_total, H HOGA_total,
#             I HOGF_pg, J HOGA_pg, K AOGF_total, L AOGA_total,
#             M AOGF_pg, N AOGA_pg
# ---------------------------------------------------------------------------
$wsHA = $wb.Worksheets.Item("Shots_HA")

$asOfNew = "2025-11-27T19:30:00Z"
for ($r = 2; $r -le 23; $r++) {
  $wsHA.Cells.Item($r, 4).Value = $asOfNew
}

$haUpdates = @(
  @{row=2;  E=14; F=15; G=449; H=404; I=32.1; J=28.9; K=531; L=506; M=35.4; N=33.7},
  @{row=3;  E=12; F=18; G=354; H=361; I=29.5; J=30.1; K=518; L=566; M=28.8; N=31.4},
  @{row=7;  E=19; F=12; G=608; H=595; I=32;   J=31.3; K=348; L=402; M=29;   N=33.5},
  @{row=11; E=13; F=17; G=353; H=467; I=27.2; J=35.9; K=404; L=648; M=23.8; N=38.1},
  @{row=13; E=16; F=14; G=607; H=398; I=37.9; J=24.9; K=411; L=415; M=29.4; N=29.6},
  @{row=15;                                        K=385;              M=29.6},
  @{row=16; E=12; F=18; G=327; H=347; I=27.2; J=28.9; K=512; L=539; M=28.4; N=29.9},
  @{row=17;                         H=344},
  @{row=19; E=18; F=12; G=606; H=511; I=33.7; J=28.4; K=393; L=408; M=32.8; N=34},
  @{row=23; E=14; F=16; G=351; H=386; I=25.1; J=27.6; K=414; L=425; M=25.9; N=26.6}
)

$haCols = @{E=5; F=6; G=7; H=8; I=9; J=10; K=11; L=12; M=13; N=14}

foreach ($u in $haUpdates) {
  $r = $u.row
  foreach ($col in $haCols.Keys) {
    if ($u.ContainsKey($col)) {
      $wsHA.Cells.Item($r, $haCols[$col]).Value = $u[$col]
    }
  }
}

# ---------------------------------------------------------------------------
# 3) Shots_Summary: as_of_utc bumped for every team, plus SOG/SOGA totals
#    for teams who played on the 27th.
#    Columns: A Team, C season_id, D as_of_utc,
#             E GP_total, F SOG_total, G SOGA_total, H SOG_pg, I SOGA_pg
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Shots_Summary")

for ($r = 2; $r -le 23; $r++) {
  $wsSummary.Cells.Item($r, 4).Value = $asOfNew
}

$summaryUpdates = @(
  @{row=2;  E=29; F=980;  G=910;  H=33.8; I=31.4},
  @{row=3;  E=30; F=872;  G=927;  H=29.1; I=30.9},
  @{row=7;  E=31; F=956;  G=997;           I=32.2},
  @{row=11; E=30; F=757;  G=1115; H=25.2; I=37.2},
  @{row=13; E=30; F=1018; G=813;  H=33.9; I=27.1},
  @{row=15;       F=873;           H=31.2},
  @{row=16; E=30; F=839;  G=886;  H=28;   I=29.5},
  @{row=17;                G=746;         I=24.9},
  @{row=19; E=30; F=999;  G=919;  H=33.3; I=30.6},
  @{row=23; E=30; F=765;  G=811;  H=25.5; I=27}
)

$summaryCols = @{E=5; F=6; G=7; H=8; I=9}

foreach ($u in $summaryUpdates) {
  $r = $u.row
  foreach ($col in $summaryCols.Keys) {
    if ($u.ContainsKey($col)) {
      $wsSummary.Cells.Item($r, $summaryCols[$col]).Value = $u[$col]
    }
  }
}

# ---------------------------------------------------------------------------
# 4) Meta_ext: as_of_utc bumped, build_version incremented 16 -> 17.
# ---------------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Meta_ext")
$wsMeta.Cells.Item(2, 2).Value = $asOfNew
$wsMeta.Cells.Item(2, 4).Value = 17
